$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the paragraph whose text contains a given substring.
# ---------------------------------------------------------------------------
function Find-ParaContaining($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Change 1: "Background images now fill up the entire draw area" paragraph
# becomes three paragraphs:
#   - "Glitch removed: Drawings scale to canvas, and no longer to the window size"
#   - "Drawing area can now be bound to a certain ratio"
#   - "Glitch removed: " + (_GoBack bookmark) + "Background images now fill up the entire draw area"
# ---------------------------------------------------------------------------
$pBg = Find-ParaContaining("Background images now fill up the entire draw area")
$bgStart = $pBg.Range.Start

# Detach the _GoBack bookmark before we start slicing paragraphs so it
# doesn't get dragged around by the paragraph split.
$hadGoBack = $d.Bookmarks.Exists("_GoBack")
if ($hadGoBack) {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
}

$insA = $d.Range($bgStart, $bgStart)
$insA.InsertParagraphBefore()
$insB = $d.Range($bgStart, $bgStart)
$insB.InsertParagraphBefore()

$pNew1 = $d.Paragraphs($pBg.Index - 2)
$pNew1.Range.Text = "Glitch removed: Drawings scale to canvas, and no longer to the window size"

$pNew2 = $d.Paragraphs($pBg.Index - 1)
$pNew2.Range.Text = "Drawing area can now be bound to a certain ratio"

$pBgAgain = Find-ParaContaining("Background images now fill up the entire draw area")
$bgAgainStart = $pBgAgain.Range.Start
$prefix = "Glitch removed: "
$insPrefix = $d.Range($bgAgainStart, $bgAgainStart)
$insPrefix.InsertBefore($prefix)

if ($hadGoBack) {
    $bmPos = $bgAgainStart + $prefix.Length
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------------
# Change 2: "New scroll bars have been added into the ..." -> "New slider bars ..."
# ---------------------------------------------------------------------------
$pScroll = Find-ParaContaining("bars have been added into the")
$pScroll.Range.Find.Execute("scroll", $true, $false, $false, $false, $false, $true, 1, $false, "slider", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: "Inputs now work around..." gets a "Glitch removed: " prefix.
# ---------------------------------------------------------------------------
$pInputs = Find-ParaContaining("Inputs now work around")
$insInputs = $d.Range($pInputs.Range.Start, $pInputs.Range.Start)
$insInputs.InsertBefore("Glitch removed: ")

# ---------------------------------------------------------------------------
# Change 4: "drawing scaling: ..." -> "Drawing scaling: ..." (capitalise).
# ---------------------------------------------------------------------------
$pDrawScale = Find-ParaContaining("rawing scaling: drawings will now scale")
$firstChar = $d.Range($pDrawScale.Range.Start, $pDrawScale.Range.Start + 1)
$firstChar.Text = "D"

Write-Output "done"
